$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$defaultStyle = $ws.Range("B2").Style

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "41.959.76"
$ws.Range("D2").Style = $defaultStyle
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +4.65%  "
$ws.Range("E2").Style = $defaultStyle
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.235.14"
$ws.Range("D3").Style = $defaultStyle
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +1.61%  "
$ws.Range("E3").Style = $defaultStyle
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.12%  "
$ws.Range("E4").Style = $defaultStyle
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "232.81"
$ws.Range("D5").Style = $defaultStyle
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +1.95%  "
$ws.Range("E5").Style = $defaultStyle
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.628"
$ws.Range("D6").Style = $defaultStyle
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -0.28%  "
$ws.Range("E6").Style = $defaultStyle
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "61.41"
$ws.Range("D7").Style = $defaultStyle
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -3.73%  "
$ws.Range("E7").Style = $defaultStyle
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.00%  "
$ws.Range("E8").Style = $defaultStyle
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +2.20%  "
$ws.Range("E9").Style = $defaultStyle
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "59.01"
$ws.Range("D10").Style = $defaultStyle
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +1.13%  "
$ws.Range("E10").Style = $defaultStyle
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0910"
$ws.Range("D11").Style = $defaultStyle
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +5.62%  "
$ws.Range("E11").Style = $defaultStyle
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +0.25%  "
$ws.Range("E12").Style = $defaultStyle
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "2.565.59"
$ws.Range("D13").Style = $defaultStyle
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +1.34%  "
$ws.Range("E13").Style = $defaultStyle
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "15.72"
$ws.Range("D14").Style = $defaultStyle
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -1.20%  "
$ws.Range("E14").Style = $defaultStyle
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "22.47"
$ws.Range("D15").Style = $defaultStyle
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +1.48%  "
$ws.Range("E15").Style = $defaultStyle
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.806"
$ws.Range("D16").Style = $defaultStyle
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -1.53%  "
$ws.Range("E16").Style = $defaultStyle
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "5.62"
$ws.Range("D17").Style = $defaultStyle
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +0.73%  "
$ws.Range("E17").Style = $defaultStyle
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.247.61"
$ws.Range("D18").Style = $defaultStyle
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +1.63%  "
$ws.Range("E18").Style = $defaultStyle
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "41.785.74"
$ws.Range("D19").Style = $defaultStyle
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +4.41%  "
$ws.Range("E19").Style = $defaultStyle
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0911"
$ws.Range("D20").Style = $defaultStyle
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +0.54%  "
$ws.Range("E20").Style = $defaultStyle
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "72.47"
$ws.Range("D21").Style = $defaultStyle
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +0.03%  "
$ws.Range("E21").Style = $defaultStyle
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.04"
$ws.Range("D22").Style = $defaultStyle
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -0.87%  "
$ws.Range("E22").Style = $defaultStyle
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "251.77"
$ws.Range("D23").Style = $defaultStyle
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +8.05%  "
$ws.Range("E23").Style = $defaultStyle
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +0.00%  "
$ws.Range("E24").Style = $defaultStyle
$ws.Range("B25").Value = "PancakeSwap"
$ws.Range("C25").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.40"
$ws.Range("D25").Style = $defaultStyle
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +1.84%  "
$ws.Range("E25").Style = $defaultStyle
$ws.Range("B26").Value = "Toncoin"
$ws.Range("C26").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.37"
$ws.Range("D26").Style = $defaultStyle
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -1.14%  "
$ws.Range("E26").Style = $defaultStyle
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.73"
$ws.Range("D27").Style = $defaultStyle
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +0.21%  "
$ws.Range("E27").Style = $defaultStyle
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.143"
$ws.Range("D28").Style = $defaultStyle
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +2.19%  "
$ws.Range("E28").Style = $defaultStyle
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "169.15"
$ws.Range("D29").Style = $defaultStyle
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -1.64%  "
$ws.Range("E29").Style = $defaultStyle
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "20.07"
$ws.Range("D30").Style = $defaultStyle
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -0.17%  "
$ws.Range("E30").Style = $defaultStyle
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -1.86%  "
$ws.Range("E31").Style = $defaultStyle
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -0.30%  "
$ws.Range("E32").Style = $defaultStyle
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -0.58%  "
$ws.Range("E33").Style = $defaultStyle
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.08"
$ws.Range("D34").Style = $defaultStyle
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +7.28%  "
$ws.Range("E34").Style = $defaultStyle
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +2.16%  "
$ws.Range("E35").Style = $defaultStyle
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +2.14%  "
$ws.Range("E36").Style = $defaultStyle
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -5.55%  "
$ws.Range("E37").Style = $defaultStyle
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.75"
$ws.Range("D38").Style = $defaultStyle
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -4.30%  "
$ws.Range("E38").Style = $defaultStyle
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.37"
$ws.Range("D39").Style = $defaultStyle
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -3.07%  "
$ws.Range("E39").Style = $defaultStyle
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.000258"
$ws.Range("D40").Style = $defaultStyle
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +31.60%  "
$ws.Range("E40").Style = $defaultStyle
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -0.09%  "
$ws.Range("E41").Style = $defaultStyle
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +5.17%  "
$ws.Range("E42").Style = $defaultStyle
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "4.68"
$ws.Range("D43").Style = $defaultStyle
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -6.88%  "
$ws.Range("E43").Style = $defaultStyle
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.56"
$ws.Range("D44").Style = $defaultStyle
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +2.98%  "
$ws.Range("E44").Style = $defaultStyle
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +0.20%  "
$ws.Range("E45").Style = $defaultStyle
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "99.52"
$ws.Range("D46").Style = $defaultStyle
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -3.95%  "
$ws.Range("E46").Style = $defaultStyle
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0961"
$ws.Range("D47").Style = $defaultStyle
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +3.40%  "
$ws.Range("E47").Style = $defaultStyle
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.482.93"
$ws.Range("D48").Style = $defaultStyle
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -2.41%  "
$ws.Range("E48").Style = $defaultStyle
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "16.64"
$ws.Range("D49").Style = $defaultStyle
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -4.77%  "
$ws.Range("E49").Style = $defaultStyle
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.80"
$ws.Range("D50").Style = $defaultStyle
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +0.04%  "
$ws.Range("E50").Style = $defaultStyle
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "52.96"
$ws.Range("D51").Style = $defaultStyle
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +6.15%  "
$ws.Range("E51").Style = $defaultStyle
